$d = $word.ActiveDocument

function Set-ParaXml {
    param($targetRange, $innerXml)
    $prefix = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
    $suffix = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $full = $prefix + $innerXml + $suffix
    $targetRange.InsertXML($full)
}

function Find-ParaByExactText {
    param($doc, $exactText)
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $cand = $doc.Paragraphs.Item($i)
        if ($cand.Range.Text -eq $exactText) {
            return $cand
        }
    }
    return $null
}

# --- Step 1: drop stray paragraph-mark <w:rFonts w:hint="eastAsia"/> run props ---
$target = Find-ParaByExactText $d "`t`t외부 클래스를 먼저 만든 후 내부 클래스 생성`r"
if ($target -eq $null) {
    Write-Output "MISSING TARGET (fix)"
} else {
    Set-ParaXml $target.Range '<w:p><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>외부 클래스를 먼저 만든 후 내부 클래스 생성</w:t></w:r></w:p>'
}
$target = Find-ParaByExactText $d "`t`t외부 클래스와 무관하게 생성`r"
if ($target -eq $null) {
    Write-Output "MISSING TARGET (fix)"
} else {
    Set-ParaXml $target.Range '<w:p><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>외부 클래스와 무관하게 생성</w:t></w:r></w:p>'
}
$target = Find-ParaByExactText $d "`t`t메서드를 호출할 때 생성`r"
if ($target -eq $null) {
    Write-Output "MISSING TARGET (fix)"
} else {
    Set-ParaXml $target.Range '<w:p><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>메서드를 호출할 때 생성</w:t></w:r></w:p>'
}
$target = Find-ParaByExactText $d "`t내부적으로 익명 객체가 생성되는것이다 .`r"
if ($target -eq $null) {
    Write-Output "MISSING TARGET (fix)"
} else {
    Set-ParaXml $target.Range '<w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="800"/><w:tab w:val="left" w:pos="1600"/><w:tab w:val="left" w:pos="2400"/><w:tab w:val="left" w:pos="3200"/><w:tab w:val="left" w:pos="4189"/></w:tabs></w:pPr><w:r><w:tab/></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">내부적으로 익명 객체가 </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>생성되는것이다</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>.</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>'
}
$target = Find-ParaByExactText $d "`t`t매개 변수만을 사용하도록 만든 함수로 외부 자료에 부수적인 영향이 발생하지 않도록 함`r"
if ($target -eq $null) {
    Write-Output "MISSING TARGET (fix)"
} else {
    Set-ParaXml $target.Range '<w:p><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">매개 변수만을 사용하도록 만든 </w:t></w:r><w:r><w:t xml:space="preserve">함수로 </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>외부 자료에 부수적인 영향이 발생하지 않도록 함</w:t></w:r></w:p>'
}

# --- Step 2: insert the new stream-notes section after the lambda-summary paragraph ---
$anchor1 = Find-ParaByExactText $d "`t함수를 변수처럼 사용하는 람다식이다.`r"
if ($anchor1 -eq $null) {
    Write-Output "MISSING ANCHOR1"
} else {
    $anchor1.Range.InsertParagraphAfter()
    $newRange = $anchor1.Next().Range
    Set-ParaXml $newRange '<w:p><w:r><w:t>2020-11-10</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>스트림</w:t></w:r></w:p><w:p><w:r><w:tab/></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">자료의 대상과 관계없이 동일한 연산을 수행할 수 있는 기능 </w:t></w:r><w:r><w:t>(</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>자료의 추상화)</w:t></w:r></w:p><w:p><w:r><w:tab/></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>한번 생성하고 사용한 스트림은 재사용할 수 없음</w:t></w:r></w:p><w:p><w:r><w:tab/></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>스트림 연산은 기존 자료를 변경하지 않음</w:t></w:r></w:p><w:p><w:r><w:tab/></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>중간연산과 최종 연산으로 구분된다.</w:t></w:r></w:p><w:p><w:r><w:tab/></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>최종 연산이 수행되어야 모든 연산이 적용되는 지연 연산</w:t></w:r></w:p><w:p><w:r><w:tab/></w:r></w:p><w:p><w:r><w:tab/></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">중간 연산 </w:t></w:r><w:r><w:t>– filter(</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>),map</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>()</w:t></w:r></w:p><w:p><w:r><w:tab/></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">최종연산 </w:t></w:r><w:r><w:t xml:space="preserve">– </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>forEach</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>) ,</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> count() , sum()</w:t></w:r></w:p><w:p/><w:p><w:r><w:tab/></w:r><w:r><w:t>R</w:t></w:r><w:r><w:t>edu</w:t></w:r><w:r><w:t>ce()</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">연산 </w:t></w:r></w:p><w:p><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>정의된 연산이 아닌 프로그래머가 직접 지정하는 연산을 적용</w:t></w:r></w:p><w:p><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">최종 연산으로 스트림의 요소를 소모하며 연산 수행 </w:t></w:r></w:p>'
}

# --- Step 3: append a trailing tab-only paragraph after the final (_GoBack) paragraph ---
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newLastRange = $d.Paragraphs.Item($d.Paragraphs.Count).Range
Set-ParaXml $newLastRange '<w:p><w:r><w:tab/></w:r></w:p>'

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
